$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E and F (ranking / colocação %) for rows 2-7 were stored as
# fractions (e.g. 0.82) but should now be stored as percentages (82.02),
# keeping the same "0.00%" number format on the cells.
for ($row = 2; $row -le 7; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = $cell.Value2 * 100
    }
}
